# Update dSF (column F) values for a set of rows, per repulled data / mean
# calculation fixes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 1
    4  = -2
    7  = -5
    8  = 0
    12 = -3
    18 = -9
    26 = -2
    29 = -1
    33 = -4
    34 = -5
    35 = -3
    36 = 5
    37 = 0
    41 = -4
    43 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
